# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback: the status text moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language detail sheets get
# their "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (including a hyperlink for the target
# file), and a couple of columns are widened so the new, longer text
# fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This label shows up on the Overview sheet (per-language Status
#    columns E/F) as well as on each language sheet's own Status column C.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Helper: fill in the "Latest Target File" (I), "Latest Handback
#    File" (J) and "Latest Handback DateTime" (K) columns for a given
#    language sheet/row, including turning column I into a hyperlink
#    that matches the existing "Source File Name" (A) hyperlink style.
# ---------------------------------------------------------------------
function Set-HandbackRow($ws, $targetCell, $handbackFileCell, $handbackDateCell, $targetFileName, $targetUrl, $handbackFileName, $handbackDateTime, $sourceStyleRange) {

    $rngTarget = $ws.Range($targetCell)
    $rngTarget.Value = $targetFileName

    # Reuse the exact same look (font/underline/color) as the existing
    # "Source File Name" hyperlink cells instead of letting Hyperlinks.Add
    # invent a brand new style.
    $sourceStyleRange.Copy()
    $rngTarget.PasteSpecial(-4122) # xlPasteFormats

    $ws.Hyperlinks.Add($rngTarget, $targetUrl, [Type]::Missing, [Type]::Missing, $targetFileName) | Out-Null

    $ws.Range($handbackFileCell).Value = $handbackFileName
    $ws.Range($handbackDateCell).Value = $handbackDateTime
}

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80b39e7abed67b3d0e13b1677b8680cb6d45bd67/e2e/42051b4b-10d0-4a1a-9116-321f9503dc5d.md"

# ---------------------------------------------------------------------
# 3. zh-cn sheet (rows 2 and 3)
# ---------------------------------------------------------------------
Set-HandbackRow $wsZhCn "I2" "J2" "K2" `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.md" $urlA `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.ed4033c643d169a794242035dec79380a48d5df4.zh-cn.xlf" `
    "2016-08-31 13:19:50" $wsZhCn.Range("A2")

Set-HandbackRow $wsZhCn "I3" "J3" "K3" `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.md" $urlA `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.ed4033c643d169a794242035dec79380a48d5df4.zh-cn.xlf" `
    "2016-08-31 13:19:50" $wsZhCn.Range("A2")

# ---------------------------------------------------------------------
# 4. de-de sheet (rows 2 and 3)
# ---------------------------------------------------------------------
Set-HandbackRow $wsDeDe "I2" "J2" "K2" `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.md" $urlA `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.ed4033c643d169a794242035dec79380a48d5df4.de-de.xlf" `
    "2016-08-31 13:19:57" $wsDeDe.Range("A2")

Set-HandbackRow $wsDeDe "I3" "J3" "K3" `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.md" $urlA `
    "42051b4b-10d0-4a1a-9116-321f9503dc5d.ed4033c643d169a794242035dec79380a48d5df4.de-de.xlf" `
    "2016-08-31 13:19:57" $wsDeDe.Range("A2")

# ---------------------------------------------------------------------
# 5. Column widths: widen the columns that now hold the longer status /
#    file-name / hyperlink text. (ColumnWidth is expressed in
#    "characters"; the stored OOXML width ends up a few pixels wider,
#    matching how these columns already looked before this edit.)
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668   # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668   # F

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664      # J
